$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1461.875
$ws.Range("I28").Value = 835
$ws.Range("K28").Value = 835
$ws.Range("M28").Value = -350

$ws.Range("H61").Value = 374.5
$ws.Range("I61").Value = 249
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 747
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -575
$ws.Range("N61").Value = -1844

$ws.Range("H64").Value = 6333
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 6999.5
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 6999.5
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -7495.5

$ws.Range("H67").Value = 6333
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 6999.5
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 6999.5
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -8715.5

$ws.Range("H70").Value = 9214.143
$ws.Range("I70").Value = 2125
$ws.Range("K70").Value = 6375
$ws.Range("M70").Value = -6105

$ws.Range("H73").Value = 9214.143
$ws.Range("I73").Value = 2125
$ws.Range("K73").Value = 6375
$ws.Range("M73").Value = -5439

$ws.Range("H74").Value = 5858
$ws.Range("J74").Value = 7650
$ws.Range("L74").Value = 7650
$ws.Range("N74").Value = -9522

$ws.Range("H76").Value = 4606.727
$ws.Range("I76").Value = 4670.5713
$ws.Range("K76").Value = 4670.5713
$ws.Range("M76").Value = -4355.5713

$ws.Range("H77").Value = 5858
$ws.Range("J77").Value = 7650
$ws.Range("L77").Value = 38250
$ws.Range("N77").Value = -47610

$ws.Range("H79").Value = 4606.727
$ws.Range("I79").Value = 4670.5713
$ws.Range("K79").Value = 4670.5713
$ws.Range("M79").Value = -3578.5713

$ws.Range("H100").Value = 4704.0835
$ws.Range("I100").Value = 3431.125
$ws.Range("J100").Value = 7250
$ws.Range("K100").Value = 3431.125
$ws.Range("L100").Value = 7250
$ws.Range("M100").Value = -2890.125
$ws.Range("N100").Value = -8332

$ws.Range("H111").Value = 853.7143

$ws.Range("H113").Value = 13734.875
$ws.Range("I113").Value = 12612.833
$ws.Range("K113").Value = 12612.833
$ws.Range("M113").Value = -9358.833000000001

$ws.Range("H137").Value = 10161.927
$ws.Range("I137").Value = 2664.5881
$ws.Range("J137").Value = 15472.542
$ws.Range("K137").Value = 7993.7643
$ws.Range("L137").Value = 46417.626
$ws.Range("M137").Value = -5443.7643
$ws.Range("N137").Value = -51517.626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 52
$ws.Range("I4").Value = 52.5
$ws.Range("K4").Value = 52.5
$ws.Range("M4").Value = 63.5

$ws.Range("H32").Value = 11779.473
$ws.Range("I32").Value = 4310.5386
$ws.Range("K32").Value = 4310.5386
$ws.Range("M32").Value = -4023.5386

$ws.Range("H61").Value = 12880.704
$ws.Range("I61").Value = 3699.6
$ws.Range("K61").Value = 3699.6
$ws.Range("M61").Value = -3487.6

$ws.Range("H68").Value = 80099
$ws.Range("J68").Value = 80099
$ws.Range("L68").Value = 80099
$ws.Range("N68").Value = -81721

$ws.Range("H71").Value = 80099
$ws.Range("J71").Value = 80099
$ws.Range("L71").Value = 240297
$ws.Range("N71").Value = -248409

$ws.Range("H76").Value = 49000
$ws.Range("J76").Value = 49000
$ws.Range("L76").Value = 49000
$ws.Range("N76").Value = -49676

$ws.Range("H79").Value = 49000
$ws.Range("J79").Value = 49000
$ws.Range("L79").Value = 49000
$ws.Range("N79").Value = -51340

$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50722

$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52496

$ws.Range("H88").Value = 2191.875
$ws.Range("J88").Value = 2362.1428
$ws.Range("L88").Value = 2362.1428
$ws.Range("N88").Value = -3174.1428

$ws.Range("H91").Value = 2191.875
$ws.Range("J91").Value = 2362.1428
$ws.Range("L91").Value = 2362.1428
$ws.Range("N91").Value = -5170.1428

$ws.Range("H97").Value = 4436.3125
$ws.Range("I97").Value = 1340.3334
$ws.Range("J97").Value = 13724.25
$ws.Range("K97").Value = 1340.3334
$ws.Range("L97").Value = 13724.25
$ws.Range("M97").Value = -844.3334
$ws.Range("N97").Value = -14716.25

$ws.Range("H132").Value = 3583814.8
$ws.Range("I132").Value = 4689.15
$ws.Range("J132").Value = 12531628
$ws.Range("K132").Value = 14067.45
$ws.Range("L132").Value = 37594884
$ws.Range("M132").Value = -11537.45
$ws.Range("N132").Value = -37599944

$ws.Range("H136").Value = 12880.704
$ws.Range("I136").Value = 3699.6
$ws.Range("K136").Value = 11098.8
$ws.Range("M136").Value = -8548.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 4484.25
$ws.Range("I36").Value = 4484.25
$ws.Range("K36").Value = 4484.25
$ws.Range("M36").Value = -3950.25

$ws.Range("H105").Value = 2341.6667
$ws.Range("I105").Value = 1831.909
$ws.Range("K105").Value = 1831.909
$ws.Range("M105").Value = -84.90900000000011

$ws.Range("H132").Value = 93883.78
$ws.Range("J132").Value = 93883.78
$ws.Range("L132").Value = 93883.78
$ws.Range("N132").Value = -104003.78

$ws.Range("H134").Value = 14471.107
$ws.Range("I134").Value = 7850
$ws.Range("K134").Value = 23550
$ws.Range("M134").Value = -21015

$ws.Range("H137").Value = 59953
$ws.Range("J137").Value = 59953
$ws.Range("L137").Value = 59953
$ws.Range("N137").Value = -70153

$ws.Range("H139").Value = 20000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11645.913
$ws.Range("I58").Value = 7739.636
$ws.Range("J58").Value = 12873.6
$ws.Range("K58").Value = 7739.636
$ws.Range("L58").Value = 12873.6
$ws.Range("M58").Value = -7536.636
$ws.Range("N58").Value = -13279.6

$ws.Range("H97").Value = 17000
$ws.Range("J97").Value = 17000
$ws.Range("L97").Value = 17000
$ws.Range("N97").Value = -18982

$ws.Range("H99").Value = 8137.8486
$ws.Range("I99").Value = 2563.9375
$ws.Range("J99").Value = 13383.883
$ws.Range("K99").Value = 2563.9375
$ws.Range("L99").Value = 13383.883
$ws.Range("M99").Value = -1065.9375
$ws.Range("N99").Value = -16379.883

$ws.Range("H126").Value = 8137.8486
$ws.Range("I126").Value = 2563.9375
$ws.Range("J126").Value = 13383.883
$ws.Range("K126").Value = 7691.8125
$ws.Range("L126").Value = 40151.649
$ws.Range("M126").Value = -5221.8125
$ws.Range("N126").Value = -45091.649

$ws.Range("H132").Value = 8311.317999999999
$ws.Range("I132").Value = 2773.4614
$ws.Range("J132").Value = 16310.444
$ws.Range("K132").Value = 8320.3842
$ws.Range("L132").Value = 48931.33199999999
$ws.Range("M132").Value = -5790.3842
$ws.Range("N132").Value = -53991.33199999999

$ws.Range("H134").Value = 26321812
$ws.Range("I134").Value = 1760.6842
$ws.Range("K134").Value = 5282.0526
$ws.Range("M134").Value = -2747.0526

$ws.Range("H136").Value = 11645.913
$ws.Range("I136").Value = 7739.636
$ws.Range("J136").Value = 12873.6
$ws.Range("K136").Value = 23218.908
$ws.Range("L136").Value = 38620.8
$ws.Range("M136").Value = -20668.908
$ws.Range("N136").Value = -43720.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 177.92308
$ws.Range("I7").Value = 153.33333
$ws.Range("J7").Value = 281.2
$ws.Range("K7").Value = 459.99999
$ws.Range("L7").Value = 843.5999999999999
$ws.Range("M7").Value = -347.99999
$ws.Range("N7").Value = -1067.6

$ws.Range("H49").Value = 1327.7142
$ws.Range("I49").Value = 149
$ws.Range("K49").Value = 447
$ws.Range("M49").Value = -291

$ws.Range("H108").Value = 1482.1
$ws.Range("I108").Value = 665.125
$ws.Range("K108").Value = 1995.375
$ws.Range("M108").Value = 884.625

$ws.Range("H121").Value = 2114.25
$ws.Range("J121").Value = 2882.8
$ws.Range("L121").Value = 8648.400000000001
$ws.Range("N121").Value = -11268.4

$ws.Range("H131").Value = 1483.75
$ws.Range("J131").Value = 1494.3368
$ws.Range("L131").Value = 4483.0104
$ws.Range("N131").Value = -14563.0104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4520.075
$ws.Range("I132").Value = 1624.4375
$ws.Range("K132").Value = 4873.3125
$ws.Range("M132").Value = -2343.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2085.2144
$ws.Range("I16").Value = 2214.923
$ws.Range("J16").Value = 399
$ws.Range("K16").Value = 2214.923
$ws.Range("L16").Value = 399
$ws.Range("M16").Value = -2044.923
$ws.Range("N16").Value = -739
